$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header D1: unit changed from kPa to Pa (labels for B1/C1 stay the same text)
$ws.Range("D1").Value = "Stiffness (Pa)"

# Update existing row 2 data with new collected values
$ws.Range("B2").Value = 8100
$ws.Range("C2").Value = [double]"9.9999999999999992E-25"
$ws.Range("D2").Value = 2755.13

# Add new data rows 3-7 collected by Rohan
$ws.Range("A3").Value = 3
$ws.Range("B3").Value = 8100
$ws.Range("C3").Formula = "=1/(5*60)"
$ws.Range("D3").Value = 4431.4399999999996

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 8100
$ws.Range("C4").Formula = "=1/(10*60)"
$ws.Range("D4").Value = 7661.76

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 8100
$ws.Range("C5").Formula = "=1/(15*60)"
$ws.Range("D5").Value = 8860.2000000000007

$ws.Range("A6").Value = 3
$ws.Range("B6").Value = 8100
$ws.Range("C6").Formula = "=1/(20*60)"
$ws.Range("D6").Value = 10619.1

$ws.Range("A7").Value = 3
$ws.Range("B7").Value = 8100
$ws.Range("C7").Formula = "=1/(25*60)"
$ws.Range("D7").Value = 12749.44

# Update the active selection as saved in the file
$ws.Range("E24").Select()
